# Update cryptos list snapshot (GitHub Actions style refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $text) {
    # Price strings like "192.59" or "1.00" would otherwise be auto-coerced
    # to numbers by Excel's type inference. Force text entry, then drop the
    # cell back to the default "Normal" style so no stray number-format
    # style sticks around (matches the source file's plain inline strings).
    $r = $ws.Range($rangeRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.507.61"
$ws.Range("E2").Value = "  +3.11%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.369.70"
$ws.Range("E3").Value = "  +4.60%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - Solana
Set-TextValue "D5" "192.59"
$ws.Range("E5").Value = "  +5.33%  "

# Row 6 - BNB
Set-TextValue "D6" "593.14"
$ws.Range("E6").Value = "  +2.69%  "

# Row 7 - was USDC, now XRP (rows 7 and 8 swapped)
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D7" "0.609"
$ws.Range("E7").Value = "  +0.97%  "

# Row 8 - was XRP, now USDC
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +3.20%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +3.77%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.422"
$ws.Range("E11").Value = "  +2.43%  "

# Row 12 - Wrapped liquid staked Ether 2.0
$ws.Range("D12").Value = "3.958.03"
$ws.Range("E12").Value = "  +4.79%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.27%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +3.34%  "

# Row 15 - Wrapped BTC
$ws.Range("D15").Value = "69.534.45"
$ws.Range("E15").Value = "  +3.06%  "

# Row 16 - Shiba Inu
$ws.Range("E16").Value = "  +2.24%  "

# Row 17 - Wrapped Ether
$ws.Range("D17").Value = "3.367.38"
$ws.Range("E17").Value = "  +3.75%  "

# Row 18 - Bitcoin Cash
Set-TextValue "D18" "449.80"
$ws.Range("E18").Value = "  +13.76%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +1.60%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +3.30%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +3.74%  "

# Row 22 - Litecoin
Set-TextValue "D22" "73.77"
$ws.Range("E22").Value = "  +3.87%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 - Wrapped eETH
$ws.Range("D24").Value = "3.520.25"

# Row 25 - Polygon
$ws.Range("E25").Value = "  +1.08%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +4.07%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +4.68%  "

# Row 28 - Internet Computer (DFINITY)
Set-TextValue "D28" "9.59"
$ws.Range("E28").Value = "  +0.36%  "

# Row 29 - Binance-Peg BSC-USD
$ws.Range("E29").Value = "  +0.07%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +2.58%  "

# Row 31 - Ethereum Classic
Set-TextValue "D31" "23.25"
$ws.Range("E31").Value = "  +2.88%  "

# Row 32 - NEAR Protocol
Set-TextValue "D32" "5.62"
$ws.Range("E32").Value = "  +1.06%  "

# Row 33 - Fetch.AI
Set-TextValue "D33" "1.30"
$ws.Range("E33").Value = "  +4.00%  "

# Row 34 - Aptos
$ws.Range("E34").Value = "  +1.52%  "

# Row 35 - USDe
$ws.Range("E35").Value = "  +0.00%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +3.79%  "

# Row 37 - Monero
Set-TextValue "D37" "165.00"

# Row 38 - Stacks
$ws.Range("E38").Value = "  +4.14%  "

# Row 39 - EnergySwap
Set-TextValue "D39" "27.29"
$ws.Range("E39").Value = "  +4.19%  "

# Row 40 - Mantle
Set-TextValue "D40" "0.820"
$ws.Range("E40").Value = "  +2.22%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +1.40%  "

# Row 42 - RenderToken
Set-TextValue "D42" "6.54"

# Row 43 - Maker
$ws.Range("D43").Value = "2.742.57"
$ws.Range("E43").Value = "  +6.00%  "

# Row 44 - dogwifhat
$ws.Range("E44").Value = "  +3.50%  "

# Row 45 - Injective Protocol
Set-TextValue "D45" "25.66"

# Row 46 - Hedera
$ws.Range("E46").Value = "  +1.16%  "

# Row 47 - Bittensor
Set-TextValue "D47" "343.72"
$ws.Range("E47").Value = "  +3.14%  "

# Row 48 - OKB
Set-TextValue "D48" "40.82"
$ws.Range("E48").Value = "  +0.82%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +3.28%  "

# Row 50 - Arweave
Set-TextValue "D50" "33.04"
$ws.Range("E50").Value = "  +8.09%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  +7.65%  "
